$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata" (first worksheet)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value added
$ws1.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact / No display for ContactDetail" row (old row 11);
# this shifts every following row up by one, so rows 12-21 land on 11-20 already
# carrying the correct (unchanged) content.
$ws1.Rows.Item(11).Delete()

# The remaining "Contact" row (now row 10) becomes the new Jurisdiction row.
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------------
# Sheet "Elements" (second worksheet)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# Short / Definition text for the extension root element
$ws2.Range("K2").Value = "Referral Indicator"
$ws2.Range("L2").Value = "Indicates whether the service resulted from a referral"
